# Cambio ppt SOLID rubrica 2026
# Updates slide 51 (sldId 305) and slide 52 (sldId 306) rubric text, and
# refreshes the cached "datetimeFigureOut" field text (3/02/2026 -> 11/02/2026)
# on the 5 slide layouts that carry a date placeholder.

$p = $ppt.ActivePresentation

function Get-ShapeById($slide, $id) {
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $sh = $slide.Shapes.Item($i)
        if ($sh.Id -eq $id) {
            return $sh
        }
    }
    return $null
}

# ---------------------------------------------------------------------
# Slide 51 (sldId 305) - "Forma de Entrega" shape (id=6)
# ---------------------------------------------------------------------
$s305 = $p.Slides.Item(51)
$shEntrega = Get-ShapeById $s305 6
$trEntrega = $shEntrega.TextFrame.TextRange

$trEntrega.Paragraphs(2).Runs(1).Text = "Cada grupo (de máximo 3 integrantes) realizará un vídeo de máximo 15 minutos y me lo comparte (por favor no adjuntarlo). A más tardar el lunes 16 de febrero antes de media noche."

$trEntrega.Paragraphs(3).Runs(1).Text = "Deberán exponer con lujo de detalles, el rediseño aplicando los principios SOLID y la implementación de esta forma:"

# Last bullet ("5 Minutos restantes, preguntas") becomes an empty,
# demoted (indent level 2 / lvl=1) bullet.
$para7 = $trEntrega.Paragraphs(7)
$para7.Runs(1).Text = ""
$para7.IndentLevel = 2

# ---------------------------------------------------------------------
# Slide 52 (sldId 306) - rubric shapes
# ---------------------------------------------------------------------
$s306 = $p.Slides.Item(52)

# "Rectángulo 4" (id=5): point breakdown list
$shPuntos = Get-ShapeById $s306 5
$trPuntos = $shPuntos.TextFrame.TextRange

$trPuntos.Paragraphs(2).Runs(1).Text = "15 puntos: el trabajo de grupo que realizaron con el proyecto MVC entregado"
# Remove the "5 puntos: de las actividades..." bullet entirely.
$trPuntos.Paragraphs(3).Delete()
# What used to be paragraph 4 is now paragraph 3 after the delete.
$trPuntos.Paragraphs(3).Runs(1).Text = "35 puntos: Aplicación de los principios SOLID a un proyecto previo hecho aplicando el paradigma de objetos. UML notación Extendida+ Implementación en C#  otro lenguaje orientado a objetos"

# "Rectángulo 1" (id=2): "12 puntos" -> "10 puntos"
$shArg = Get-ShapeById $s306 2
$shArg.TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "10 puntos"

# "Rectángulo 6" (id=7): "18 puntos" -> "15 puntos"
$shImpl = Get-ShapeById $s306 7
$shImpl.TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "15 puntos"

# "Rectángulo 7" (id=8): "12puntos" -> "10 puntos"
$shUml = Get-ShapeById $s306 8
$shUml.TextFrame.TextRange.Paragraphs(2).Runs(1).Text = "10 puntos"

# ---------------------------------------------------------------------
# Refresh the cached date placeholder text on every slide layout that
# shows "3/02/2026" (Solo el título, Encabezado de sección, Título y
# objetos, Dos objetos, Diapositiva de título).
# ---------------------------------------------------------------------
$master = $p.SlideMaster
for ($li = 1; $li -le $master.CustomLayouts.Count; $li++) {
    $layout = $master.CustomLayouts.Item($li)
    for ($si = 1; $si -le $layout.Shapes.Count; $si++) {
        $lsh = $layout.Shapes.Item($si)
        if ($lsh.HasTextFrame) {
            $ltr = $lsh.TextFrame.TextRange
            if ($ltr.Text -eq "3/02/2026") {
                $ltr.Text = "11/02/2026"
            }
        }
    }
}
